$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 24,14
$values[0,0] = 13.50883376864503
$values[0,1] = 12.35815665676973
$values[0,2] = 0
$values[0,3] = 17.16067776260198
$values[0,4] = 34.81067705249825
$values[0,5] = 25.69717249075682
$values[0,6] = 13.53424162797267
$values[0,7] = 21.65866818436182
$values[0,8] = 7.299823904988677
$values[0,9] = 0
$values[0,10] = 12.9688422575231
$values[0,11] = 15.92527140931689
$values[0,12] = 17.84394513949599
$values[0,13] = 20.1894045527674
$values[1,0] = 13.06193415403822
$values[1,1] = 12.2820439440475
$values[1,2] = 0
$values[1,3] = 17.20220723113597
$values[1,4] = 34.85120485821498
$values[1,5] = 25.71724426280887
$values[1,6] = 13.57460850032949
$values[1,7] = 21.75191486764657
$values[1,8] = 7.281398475968065
$values[1,9] = 0
$values[1,10] = 12.96580348732661
$values[1,11] = 15.82845165467315
$values[1,12] = 17.88285341267488
$values[1,13] = 20.24752042497905
$values[2,0] = 12.78067532156622
$values[2,1] = 12.23493026959579
$values[2,2] = 0
$values[2,3] = 17.2296987749856
$values[2,4] = 34.88398574801081
$values[2,5] = 25.7381694116998
$values[2,6] = 13.60151559612721
$values[2,7] = 21.8131306562745
$values[2,8] = 7.269905538331855
$values[2,9] = 0
$values[2,10] = 12.96549481223798
$values[2,11] = 15.77052344032688
$values[2,12] = 17.90839473683354
$values[2,13] = 20.2875624634345
$values[3,0] = 12.66450285316686
$values[3,1] = 12.21564405881936
$values[3,2] = 0
$values[3,3] = 17.24140338453747
$values[3,4] = 34.89932856745526
$values[3,5] = 25.74885465837759
$values[3,6] = 13.61301392688847
$values[3,7] = 21.83907271833163
$values[3,8] = 7.265177593837375
$values[3,9] = 0
$values[3,10] = 12.96576214204716
$values[3,11] = 15.7473168106817
$values[3,12] = 17.91921905637551
$values[3,13] = 20.30497393275182
$values[4,0] = 12.64512364815934
$values[4,1] = 12.21243659347056
$values[4,2] = 0
$values[4,3] = 17.24337723816473
$values[4,4] = 34.90199602667823
$values[4,5] = 25.75075911045253
$values[4,6] = 13.61495543496108
$values[4,7] = 21.8434405385737
$values[4,8] = 7.264389855992571
$values[4,9] = 0
$values[4,10] = 12.9658303185499
$values[4,11] = 15.74348800194747
$values[4,12] = 17.92104157514764
$values[4,13] = 20.30793110847435
$values[5,0] = 12.77911465139433
$values[5,1] = 12.23467050998003
$values[5,2] = 0
$values[5,3] = 17.22985459574513
$values[5,4] = 34.8841846344313
$values[5,5] = 25.7383047860618
$values[5,6] = 13.60166850669492
$values[5,7] = 21.81347648685171
$values[5,8] = 7.269841954802766
$values[5,9] = 0
$values[5,10] = 12.96549682384777
$values[5,11] = 15.77020882699273
$values[5,12] = 17.90853903194004
$values[5,13] = 20.28779285373207
$values[6,0] = 13.3562629321934
$values[6,1] = 12.33199705133613
$values[6,2] = 0
$values[6,3] = 17.17458399132472
$values[6,4] = 34.82301149891057
$values[6,5] = 25.7023057251315
$values[6,6] = 13.54771975083211
$values[6,7] = 21.68999725600446
$values[6,8] = 7.293508519633525
$values[6,9] = 0
$values[6,10] = 12.9674722741784
$values[6,11] = 15.89158324065597
$values[6,12] = 17.85701839061201
$values[6,13] = 20.20853709845094
$values[7,0] = 14.4269779638179
$values[7,1] = 12.51946245874934
$values[7,2] = 0
$values[7,3] = 17.0819792143693
$values[7,4] = 34.76573895320827
$values[7,5] = 25.70009777052061
$values[7,6] = 13.45876256873764
$values[7,7] = 21.47929694378577
$values[7,8] = 7.338463575158219
$values[7,9] = 0
$values[7,10] = 12.98362918491744
$values[7,11] = 16.14086382934804
$values[7,12] = 17.76905739829003
$values[7,13] = 20.08778578869212
$values[8,0] = 15.16853337769867
$values[8,1] = 12.65462676334002
$values[8,2] = 0
$values[8,3] = 17.02352555546686
$values[8,4] = 34.76187766127521
$values[8,5] = 25.74025522160742
$values[8,6] = 13.40367140692313
$values[8,7] = 21.34367115665822
$values[8,8] = 7.370566831144679
$values[8,9] = 0
$values[8,10] = 13.00288563076741
$values[8,11] = 16.32979085674891
$values[8,12] = 17.71235591248038
$values[8,13] = 20.02031393444286
$values[9,0] = 15.49470643054423
$values[9,1] = 12.71545352456103
$values[9,2] = 0
$values[9,3] = 16.99900641917653
$values[9,4] = 34.76840301587701
$values[9,5] = 25.76757992798031
$values[9,6] = 13.38083840382118
$values[9,7] = 21.28613768783227
$values[9,8] = 7.38496232737272
$values[9,9] = 0
$values[9,10] = 13.01322542205912
$values[9,11] = 16.41675508501957
$values[9,12] = 17.68827214254797
$values[9,13] = 19.99425246714944
$values[10,0] = 15.61651581073384
$values[10,1] = 12.7383831524322
$values[10,2] = 0
$values[10,3] = 16.99001896089966
$values[10,4] = 34.77206200738737
$values[10,5] = 25.77922520872951
$values[10,6] = 13.37251259799824
$values[10,7] = 21.26495032742755
$values[10,8] = 7.39038280028179
$values[10,9] = 0
$values[10,10] = 13.017365725376
$values[10,11] = 16.44981237060573
$values[10,12] = 17.6793974522504
$values[10,13] = 19.98505110217704
$values[11,0] = 15.59035928974142
$values[11,1] = 12.73344962473109
$values[11,2] = 0
$values[11,3] = 16.99194135168143
$values[11,4] = 34.77122119854661
$values[11,5] = 25.77665954821263
$values[11,6] = 13.37429145032887
$values[11,7] = 21.26948673952734
$values[11,8] = 7.389216786876984
$values[11,9] = 0
$values[11,10] = 13.0164640754791
$values[11,11] = 16.44268763119969
$values[11,12] = 17.68129787456188
$values[11,13] = 19.98700306956143
$values[12,0] = 15.50476246821432
$values[12,1] = 12.71734209863293
$values[12,2] = 0
$values[12,3] = 16.99826105800657
$values[12,4] = 34.76868025438174
$values[12,5] = 25.76851200950805
$values[12,6] = 13.38014700958039
$values[12,7] = 21.28438257605506
$values[12,8] = 7.385408894700027
$values[12,9] = 0
$values[12,10] = 13.01356155489018
$values[12,11] = 16.41947232644566
$values[12,12] = 17.68753710336657
$values[12,13] = 19.99348207876176
$values[13,0] = 15.4521071028616
$values[13,1] = 12.7074619385327
$values[13,2] = 0
$values[13,3] = 17.00217077668663
$values[13,4] = 34.76727846058348
$values[13,5] = 25.76369028601982
$values[13,6] = 13.3837754585587
$values[13,7] = 21.29358477378574
$values[13,8] = 7.383072416801546
$values[13,9] = 0
$values[13,10] = 13.01181288737269
$values[13,11] = 16.40526804440582
$values[13,12] = 17.69139074062588
$values[13,13] = 19.997537634581
$values[14,0] = 15.14698374018898
$values[14,1] = 12.65063751796341
$values[14,2] = 0
$values[14,3] = 17.02516958397489
$values[14,4] = 34.76161763511632
$values[14,5] = 25.73865140306347
$values[14,6] = 13.40520844742427
$values[14,7] = 21.34751493480975
$values[14,8] = 7.369621787728341
$values[14,9] = 0
$values[14,10] = 13.00224149604654
$values[14,11] = 16.3241262886347
$values[14,12] = 17.7139642028878
$values[14,13] = 20.02211044884623
$values[15,0] = 14.95686615347846
$values[15,1] = 12.61560191567823
$values[15,2] = 0
$values[15,3] = 17.0398088672498
$values[15,4] = 34.76026427502805
$values[15,5] = 25.72560800372277
$values[15,6] = 13.41892770065264
$values[15,7] = 21.38166611870059
$values[15,8] = 7.36131629846108
$values[15,9] = 0
$values[15,10] = 12.99677279080856
$values[15,11] = 16.27459496745387
$values[15,12] = 17.72824984538584
$values[15,13] = 20.03837251780254
$values[16,0] = 14.84646989993783
$values[16,1] = 12.59538877427103
$values[16,2] = 0
$values[16,3] = 17.0484240268184
$values[16,4] = 34.7602656191593
$values[16,5] = 25.71895863542617
$values[16,6] = 13.42702837747652
$values[16,7] = 21.40170081401468
$values[16,8] = 7.356519631002318
$values[16,9] = 0
$values[16,10] = 12.99377619595276
$values[16,11] = 16.24620316754034
$values[16,12] = 17.73662755863191
$values[16,13] = 20.04816192024679
$values[17,0] = 14.80891527318288
$values[17,1] = 12.58853465573112
$values[17,2] = 0
$values[17,3] = 17.05137448394526
$values[17,4] = 34.76040007420968
$values[17,5] = 25.71685385523968
$values[17,6] = 13.42980714355907
$values[17,7] = 21.40855149757166
$values[17,8] = 7.354892225462048
$values[17,9] = 0
$values[17,10] = 12.99278723652363
$values[17,11] = 16.23660753881501
$values[17,12] = 17.73949177663589
$values[17,13] = 20.05155125969492
$values[18,0] = 14.97721347656447
$values[18,1] = 12.61933795662772
$values[18,2] = 0
$values[18,3] = 17.03823030923304
$values[18,4] = 34.76032765869495
$values[18,5] = 25.7269082550289
$values[18,6] = 13.41744555590647
$values[18,7] = 21.37799010871938
$values[18,8] = 7.362202463252761
$values[18,9] = 0
$values[18,10] = 12.9973395539616
$values[18,11] = 16.27985774132236
$values[18,12] = 17.72671245684336
$values[18,13] = 20.03659626833191
$values[19,0] = 15.52995134809223
$values[19,1] = 12.7220761691011
$values[19,2] = 0
$values[19,3] = 16.99639673957689
$values[19,4] = 34.76939437762464
$values[19,5] = 25.77086995511099
$values[19,6] = 13.37841838902555
$values[19,7] = 21.27999103737391
$values[19,8] = 7.386528207082479
$values[19,9] = 0
$values[19,10] = 13.01440801254169
$values[19,11] = 16.42628797681332
$values[19,12] = 17.68569783758594
$values[19,13] = 19.99156090865864
$values[20,0] = 15.8812196063071
$values[20,1] = 12.78861026374432
$values[20,2] = 0
$values[20,3] = 16.97078936025414
$values[20,4] = 34.78224247841288
$values[20,5] = 25.80716397621599
$values[20,6] = 13.35478038539082
$values[20,7] = 21.2194365636219
$values[20,8] = 7.402246641242876
$values[20,9] = 0
$values[20,10] = 13.02687264872832
$values[20,11] = 16.52271219279981
$values[20,12] = 17.66032200773162
$values[20,13] = 19.96601923486456
$values[21,0] = 15.69468410967027
$values[21,1] = 12.75315872159259
$values[21,2] = 0
$values[21,3] = 16.98429807304002
$values[21,4] = 34.7747529983548
$values[21,5] = 25.78710310215245
$values[21,6] = 13.36722542200429
$values[21,7] = 21.25143574699486
$values[21,8] = 7.393874144377952
$values[21,9] = 0
$values[21,10] = 13.0201010345023
$values[21,11] = 16.47118950535046
$values[21,12] = 17.67373494860117
$values[21,13] = 19.97929478465079
$values[22,0] = 14.96801785513698
$values[22,1] = 12.61764911206281
$values[22,2] = 0
$values[22,3] = 17.03894335591289
$values[22,4] = 34.76029657494735
$values[22,5] = 25.72631776512336
$values[22,6] = 13.41811496915312
$values[22,7] = 21.37965078464723
$values[22,8] = 7.361801895776914
$values[22,9] = 0
$values[22,10] = 12.99708286091815
$values[22,11] = 16.27747817671307
$values[22,12] = 17.72740699710634
$values[22,13] = 20.03739793988199
$values[23,0] = 14.14472633733535
$values[23,1] = 12.46916270655547
$values[23,2] = 0
$values[23,3] = 17.10534567641729
$values[23,4] = 34.77451821680793
$values[23,5] = 25.69335797387386
$values[23,6] = 13.48102518523998
$values[23,7] = 21.53293019728509
$values[23,8] = 7.326461719972492
$values[23,9] = 0
$values[23,10] = 12.97795335533366
$values[23,11] = 16.07233143282579
$values[23,12] = 17.7914585453055
$values[23,13] = 20.11672879338586

$ws.Range("B2:O25").Value2 = $values
Write-Host "Updated loading_percent values for rows 2-25, columns B-O"